$d = $word.ActiveDocument
$d.Content.Find.Execute("test", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Test 23", 2)
